# Updates the North Macedonia 1.MFL 2023-2024 results sheet:
# - Rows 67-76 are re-shuffled (a previously-missing fixture, Shkendija vs
#   Voska Sport, is reinserted into its correct chronological slot and the
#   surrounding rows shift to match), and
# - Two newly played fixtures (rows 77-78) are appended.
# Columns: A=Indice, B=pais, C=torneio, D=temporada, E=data_partida,
# F=home, G=home_ft_gols, H=away, I=away_ft_gols, J..U=odds/timestamps, V=url

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-write columns F:V for rows 67-76 with the corrected match order ---
# Each entry: row, home, home_goals, away, away_goals,
#             home_open_odds, home_open_dt, home_close_odds, home_close_dt,
#             draw_open_odds, draw_open_dt, draw_close_odds, draw_close_dt,
#             away_open_odds, away_open_dt, away_close_odds, away_close_dt, url
$rowsData = @(
        @(67, 'Shkendija', 1, 'Voska Sport', 1, 1.34, '28/10/2023 02:13', 1.32, '29/10/2023 00:12', 4.11, '28/10/2023 02:13', 4.53, '29/10/2023 12:38', 6.67, '28/10/2023 02:13', 8.06, '29/10/2023 12:38', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/shkendija-tetovo-voska-sport/808PO43b/'),
        @(68, 'Vardar', 1, 'Brera Strumica', 0, 3.01, '28/10/2023 02:13', 2.86, '29/10/2023 12:56', 2.67, '28/10/2023 02:13', 2.86, '29/10/2023 12:56', 2.27, '28/10/2023 02:13', 2.47, '29/10/2023 12:56', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/vardar-brera-strumica/U99LPplh/'),
        @(69, 'Makedonija GP', 1, 'Tikves', 1, 2.08, '28/10/2023 13:43', 2.18, '29/10/2023 12:41', 2.85, '28/10/2023 13:43', 2.79, '29/10/2023 12:41', 3.16, '28/10/2023 13:43', 3.48, '29/10/2023 12:41', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/makedonija-gp-tikves/2a1YMrYA/'),
        @(70, 'KF Gostivar', 2, 'Struga', 1, 3.56, '28/10/2023 03:12', 3.26, '29/10/2023 12:57', 2.94, '28/10/2023 03:12', 2.96, '29/10/2023 12:34', 1.89, '28/10/2023 03:12', 2.17, '29/10/2023 12:57', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/kf-gostivar-struga/jPcUNOI4/'),
        @(71, 'Shkupi', 3, 'Rabotnicki', 0, 1.47, '28/10/2023 02:13', 1.53, '29/10/2023 12:44', 3.61, '28/10/2023 02:13', 3.58, '29/10/2023 12:44', 5.43, '28/10/2023 02:13', 5.76, '29/10/2023 12:44', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/shkupi-rabotnicki/YoO54zIq/'),
        @(72, 'Struga', 4, 'Vardar', 0, 1.33, '31/10/2023 01:12', 1.08, '01/11/2023 12:58', 3.91, '31/10/2023 01:12', 7.98, '01/11/2023 12:59', 6.9, '31/10/2023 01:12', 27.05, '01/11/2023 12:59', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/struga-vardar/QwlA9Dtc/'),
        @(73, 'Tikves', 2, 'KF Gostivar', 2, 2.19, '03/11/2023 01:13', 2.44, '04/11/2023 12:41', 2.84, '03/11/2023 01:13', 2.61, '04/11/2023 12:41', 2.95, '03/11/2023 01:13', 3.2, '04/11/2023 12:41', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/tikves-kf-gostivar/dK90QHv3/'),
        @(74, 'Bregalnica Stip', 2, 'Makedonija GP', 3, 2.01, '03/11/2023 01:13', 2.24, '04/11/2023 12:51', 2.86, '03/11/2023 01:13', 2.89, '04/11/2023 12:51', 3.32, '03/11/2023 01:13', 3.2, '04/11/2023 12:51', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/bregalnica-stip-makedonija-gp/rH5dRcPd/'),
        @(75, 'Rabotnicki', 1, 'Sileks', 0, 2.11, '03/11/2023 01:13', 2.34, '04/11/2023 12:54', 2.83, '03/11/2023 01:13', 2.77, '04/11/2023 12:54', 3.13, '03/11/2023 01:13', 3.16, '04/11/2023 12:54', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/rabotnicki-sileks/UDHE2End/'),
        @(76, 'Voska Sport', 1, 'Brera Strumica', 1, 2.51, '03/11/2023 01:13', 2.05, '04/11/2023 12:53', 2.77, '03/11/2023 01:13', 3.03, '04/11/2023 12:53', 2.58, '03/11/2023 01:13', 3.48, '04/11/2023 12:53', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/voska-sport-brera-strumica/IZc9Oe9F/')
)

$colsFV = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

foreach ($row in $rowsData) {
    $r = $row[0]
    for ($i = 0; $i -lt $colsFV.Length; $i++) {
        $ws.Range($colsFV[$i] + $r).Value = $row[$i + 1]
    }
}

# --- Append two brand-new fixture rows (77-78), copying formatting from
#     the last existing data row (76) so styles/number formats match ---
$ws.Range("A76:V76").Copy()
$ws.Range("A77:V77").PasteSpecial(-4122)
$ws.Range("A76:V76").Copy()
$ws.Range("A78:V78").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Each entry: row, Indice, pais, torneio, temporada, data_partida,
#             home, home_goals, away, away_goals,
#             home_open_odds, home_open_dt, home_close_odds, home_close_dt,
#             draw_open_odds, draw_open_dt, draw_close_odds, draw_close_dt,
#             away_open_odds, away_open_dt, away_close_odds, away_close_dt, url
$newRowsData = @(
        @(77, 76, 'north-macedonia', '1-mfl', '2023-2024', 45235.54166666666, 'Shkupi', 2, 'Vardar', 0, 1.24, '04/11/2023 01:13', 1.17, '05/11/2023 12:57', 4.6, '04/11/2023 01:13', 5.55, '05/11/2023 12:57', 8.07, '04/11/2023 01:13', 13.91, '05/11/2023 12:57', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/shkupi-vardar/jNIA3fXk/'),
        @(78, 77, 'north-macedonia', '1-mfl', '2023-2024', 45235.54166666666, 'Struga', 2, 'Shkendija', 1, 2.17, '04/11/2023 01:13', 2.19, '05/11/2023 12:58', 2.72, '04/11/2023 01:13', 3.03, '05/11/2023 12:58', 3.13, '04/11/2023 01:13', 3.14, '05/11/2023 12:58', 'https://www.betexplorer.com/football/north-macedonia/1-mfl/struga-shkendija-tetovo/vqa5Pyf9/')
)

$colsAV = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

foreach ($row in $newRowsData) {
    $r = $row[0]
    for ($i = 0; $i -lt $colsAV.Length; $i++) {
        $ws.Range($colsAV[$i] + $r).Value = $row[$i + 1]
    }
}

"done"
